# Generate Report for Handoff
# Bumps the "Latest HO Xliff Generate Date" / "Latest Handoff Datetime" timestamps
# for the newly generated handoff files, and records their Priority ("ht").

$wb = $excel.ActiveWorkbook

$rows = @(7, 9, 11, 12, 13, 14)

# --- Overview sheet: column G "Latest HO Xliff Generate Date" ---
$wsOverview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $wsOverview.Range("G$r").Value = "2016-08-29 06:22:10"
}

# --- zh-cn sheet: column E "Priority" and column H "Latest Handoff Datetime" ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $wsZhCn.Range("E$r").Value = "ht"
    $wsZhCn.Range("H$r").Value = "2016-08-29 06:22:00"
}

# --- de-de sheet: column E "Priority" and column H "Latest Handoff Datetime" ---
$wsDeDe = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $wsDeDe.Range("E$r").Value = "ht"
    $wsDeDe.Range("H$r").Value = "2016-08-29 06:22:10"
}
